$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q1)
$ws.Range("B2").Value = -0.2425300391056354
$ws.Range("C2").Value = 0.527227898904767
$ws.Range("D2").Value = 0.9024928702735646
$ws.Range("E2").Value = 0.9499962475049912
$ws.Range("F2").Value = 0.927477568261686

# Row 3 (Q2)
$ws.Range("B3").Value = 0.09740790795838196
$ws.Range("C3").Value = 0.6911176283945465
$ws.Range("D3").Value = 1.207585491876232
$ws.Range("E3").Value = 1.098901948253907
$ws.Range("F3").Value = 1.105467835429992

# Row 4 (Q3)
$ws.Range("B4").Value = -0.004016768313543259
$ws.Range("C4").Value = 0.5629705147095596
$ws.Range("D4").Value = 0.6780094397277413
$ws.Range("E4").Value = 0.8234132861010571
$ws.Range("F4").Value = 0.8317631293628042

# Row 5 (Q4)
$ws.Range("B5").Value = 0.09599577651968134
$ws.Range("C5").Value = 0.6642539179825266
$ws.Range("D5").Value = 0.9589064666597493
$ws.Range("E5").Value = 0.9792376967109412
$ws.Range("F5").Value = 0.9846199838003967
$ws.Range("G5").Value = 49

# Row 6 (Q5)
$ws.Range("B6").Value = 0.008978782371844934
$ws.Range("C6").Value = 0.5174901747496273
$ws.Range("D6").Value = 0.5159409445670499
$ws.Range("E6").Value = 0.7182902926860768
$ws.Range("F6").Value = 0.7258347454151787
$ws.Range("G6").Value = 48

# Row 7 (Q6)
$ws.Range("B7").Value = -0.01246672580935241
$ws.Range("C7").Value = 0.4777751499603008
$ws.Range("D7").Value = 0.5121909218362894
$ws.Range("E7").Value = 0.7156751510540864
$ws.Range("F7").Value = 0.7249207688243822
$ws.Range("G7").Value = 39

# Row 8 (Q7)
$ws.Range("B8").Value = -0.008167537878776951
$ws.Range("C8").Value = 0.427281569976559
$ws.Range("D8").Value = 0.3635138582900146
$ws.Range("E8").Value = 0.6029211045319401
$ws.Range("F8").Value = 0.6109583009916096
$ws.Range("G8").Value = 38

# Row 9 (Q8)
$ws.Range("B9").Value = 0.0204133146117494
$ws.Range("C9").Value = 0.3312074429360707
$ws.Range("D9").Value = 0.1677068730083327
$ws.Range("E9").Value = 0.4095202962104964
$ws.Range("F9").Value = 0.4191117727702709
$ws.Range("G9").Value = 21

# Row 10 (Q9)
$ws.Range("B10").Value = 0.001286647701700525
$ws.Range("C10").Value = 0.3911046241481878
$ws.Range("D10").Value = 0.2110015332851695
$ws.Range("E10").Value = 0.4593490320934284
$ws.Range("F10").Value = 0.4766871486157238
$ws.Range("G10").Value = 14

# Row 11
$ws.Range("B11").Value = -0.1188039491228551
$ws.Range("C11").Value = 0.3962730645310045
$ws.Range("D11").Value = 0.2251242489150065
$ws.Range("E11").Value = 0.4744726008053642
$ws.Range("F11").Value = 0.5135779767813022
